$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-01-22 Wednesday"; new = "2025-01-23 Thursday"},
    @{old = "960÷2="; new = "152÷7="},
    @{old = "743÷9="; new = "505÷9="},
    @{old = "784÷7="; new = "992÷7="},
    @{old = "890÷8="; new = "123÷6="},
    @{old = "312÷7="; new = "219÷4="},
    @{old = "657÷4="; new = "516÷8="},
    @{old = "631÷4="; new = "234÷4="},
    @{old = "880÷5="; new = "238÷9="},
    @{old = "763÷8="; new = "667÷4="},
    @{old = "755÷4="; new = "757÷7="},
    @{old = "219÷8="; new = "283÷3="},
    @{old = "396÷4="; new = "103÷9="},
    @{old = "245÷5="; new = "733÷2="},
    @{old = "569÷8="; new = "778÷3="},
    @{old = "882÷5="; new = "865÷5="},
    @{old = "574÷3="; new = "943÷9="},
    @{old = "550÷4="; new = "587÷3="},
    @{old = "223÷4="; new = "290÷2="},
    @{old = "553÷7="; new = "461÷9="},
    @{old = "980÷9="; new = "692÷3="},
    @{old = "331÷8="; new = "106÷5="},
    @{old = "293÷3="; new = "397÷6="},
    @{old = "685÷9="; new = "839÷3="},
    @{old = "139÷5="; new = "643÷7="},
    @{old = "590÷8="; new = "778÷5="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $r.new, 2)
}
